$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsConcept = $wb.Worksheets.Item("Include from SNOMED CT")

# --- Update Metadata sheet (sheet1) ---
# Version (A3/B3)
$wsMeta.Range("B3").Value = "1.2.0"
# Date (A8/B8)
$wsMeta.Range("B8").Value = "2024-03-28T10:46:20+01:00"
# Description (A12/B12)
$wsMeta.Range("B12").Value = 'The "CH ELM Results Geni Spec" material group provides a curated set of codes representing specific materials. Each code within this group is selected to ensure consistency and accuracy for medical coding related to the primary LOINC codes. Clients using the "CH ELM Results Geni Spec" group should refer to the provided codes to ensure they capture and report material information correctly and consistently.'

# --- Update Concept sheet (sheet2) ---
$codes = @(
    @("119396006", "Specimen from endometrium (specimen)"),
    @("119324002", "Specimen of unknown material (specimen)"),
    @("119395005", "Specimen from uterine cervix (specimen)"),
    @("258527002", "Anal swab (specimen)"),
    @("122575003", "Urine specimen (specimen)"),
    @("257261003", "Swab (specimen)"),
    @("119393003", "Specimen from urethra (specimen)"),
    @("119347001", "Seminal fluid specimen (specimen)"),
    @("119394009", "Specimen from vagina (specimen)"),
    @("119344008", "Specimen from genital system (specimen)"),
    @("", ""),
    @("System URI", "http://snomed.info/sct")
)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $wsConcept.Cells.Item($row, 1).Value = $codes[$i][0]
    $wsConcept.Cells.Item($row, 2).Value = $codes[$i][1]
}
